$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("M1").Value = "i2"
$ws.Range("N1").Value = "d2"

# New data cells
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 2

$ws.Range("M3").Value = 7
$ws.Range("N3").Value = 3

$ws.Range("M4").Value = 6
$ws.Range("N4").Value = 4

$ws.Range("M5").Value = 5
$ws.Range("N5").Value = 7

$ws.Range("M6").Value = 7
$ws.Range("N6").Value = 8

# Update selection to reflect new range
$ws.Range("M1:N6").Select()
